$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new program rows: "openssh" (row 48) and "okular" (row 49) ---
# Row 43 (7-zip-like pattern: B=supported, C=not supported, D=not supported)
# is used as the format donor since its B/C/D styling (green/red/red) matches
# the target look of the two new rows.

# Row 48: openssh
$ws.Range("A48").Value = "openssh"
$ws.Range("B43:D43").Copy()
$ws.Range("B48").PasteSpecial(-4122)

# Row 49: okular
$ws.Range("A49").Value = "okular"
$ws.Range("B43:D43").Copy()
$ws.Range("B49").PasteSpecial(-4122)

# Move the active selection to where the author left it after the edit
[void]$ws.Range("E47").Select()

# --- Header / footer font name fix: "Regular" -> "Normal" ---
$ps = $ws.PageSetup
$ps.CenterHeader = '&"Times New Roman,Normal"&12&A'
$ps.CenterFooter = '&"Times New Roman,Normal"&12Página &P'
